$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Results"

# Insert 8 new rows above row 2, pushing the previous data (old rows 2:9)
# down to rows 10:17
$ws.Rows("2:9").Insert()

# --- New experiment data (rows 2-9) ---
$ws.Range("A2").Formula = "=1"
$ws.Range("B2").Value = "A->B"
$ws.Range("C2").Value = "A"
$ws.Range("D2").Value = 31.895
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 9
$ws.Range("G2").Value = 3
$ws.Range("H2").Formula = "=G2/F2"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "A->B"
$ws.Range("C3").Value = "B"
$ws.Range("D3").Value = 102.792
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 35
$ws.Range("G3").Value = 5
$ws.Range("H3").Formula = "=G3/F3"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "A->B"
$ws.Range("C4").Value = "A"
$ws.Range("D4").Value = 144.64099999999999
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 16
$ws.Range("G4").Value = 16
$ws.Range("H4").Formula = "=G4/F4"

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "A->B"
$ws.Range("C5").Value = "B"
$ws.Range("D5").Value = 35.966999999999999
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 10
$ws.Range("H5").Formula = "=G5/F5"

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "A->B"
$ws.Range("C6").Value = "A"
$ws.Range("D6").Value = 83.334000000000003
$ws.Range("E6").Value = 12
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 12
$ws.Range("H6").Formula = "=G6/F6"

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "A->B"
$ws.Range("C7").Value = "B"
$ws.Range("D7").Value = 240.61199999999999
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 10
$ws.Range("H7").Formula = "=G7/F7"

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "A->B"
$ws.Range("C8").Value = "A"
$ws.Range("D8").Value = 47.143000000000001
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 9
$ws.Range("G8").Value = 9
$ws.Range("H8").Formula = "=G8/F8"

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "A->B"
$ws.Range("C9").Value = "B"
$ws.Range("D9").Value = 64.253
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 8
$ws.Range("H9").Formula = "=G9/F9"

# --- Previous data, now shifted to rows 10-17 and updated ---
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "B->A"
$ws.Range("C10").Value = "B"
$ws.Range("D10").Value = 10.5
$ws.Range("E10").Formula = "=1"
$ws.Range("F10").Value = 10
$ws.Range("G10").Formula = "=1"
$ws.Range("H10").Formula = "=G10/F10*100"

$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "B->A"
$ws.Range("C11").Value = "A"
$ws.Range("D11").Formula = "=D10+2"
$ws.Range("E11").Formula = "=E10+1"
$ws.Range("F11").Value = 15
$ws.Range("G11").Value = 2
$ws.Range("H11").Formula = "=G11/F11*100"

$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "B->A"
$ws.Range("C12").Value = "B"
$ws.Range("D12").Value = 13.5
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 25
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 12

$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "B->A"
$ws.Range("C13").Value = "A"
$ws.Range("D13").Value = 14.5
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 30
$ws.Range("G13").Value = 4
$ws.Range("H13").Formula = "=G13/F13*100"

$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "B->A"
$ws.Range("C14").Value = "B"
$ws.Range("D14").Formula = "=D13+1"
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 53
$ws.Range("G14").Value = 2
$ws.Range("H14").Formula = "=G14/F14*100"

$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "B->A"
$ws.Range("C15").Value = "A"
$ws.Range("D15").Formula = "=D14+1"
$ws.Range("E15").Value = 6
$ws.Range("F15").Value = 7
$ws.Range("G15").Formula = "=G14+1"
$ws.Range("H15").Formula = "=G15/F15*100"

$ws.Range("A16").Value = 8
$ws.Range("B16").Value = "B->A"
$ws.Range("C16").Value = "B"
$ws.Range("D16").Formula = "=D15+1"
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 5
$ws.Range("G16").Formula = "=G15+1"
$ws.Range("H16").Value = 80

$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "B->A"
$ws.Range("C17").Value = "A"
$ws.Range("D17").Formula = "=D16+1"
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 10
$ws.Range("G17").Formula = "=G16+1"
$ws.Range("H17").Formula = "=G17/F17*100"

# --- Styling ---
# Highlight the older (shifted) data with a yellow fill
$ws.Range("A10:H17").Interior.Color = 65535

# Apply the Percent number-format style to the new Hit_Percent column
$ws.Range("H2:H9").Style = "Percent"

# Update selection to match the author's final cursor position
$ws.Range("B27").Select()
